$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-82 down to 79-83
$ws.Rows.Item(78).Insert()

# Fill in the new row 78 with its data
$ws.Cells.Item(78, 1).Value = 4
$ws.Cells.Item(78, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(78, 3).Value = "Los Lagos"
$ws.Cells.Item(78, 4).Value = 44516
$ws.Cells.Item(78, 5).Value = 10
$ws.Cells.Item(78, 6).Value = 100112022
$ws.Cells.Item(78, 7).Value = "Arveja Verde"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 20000
$ws.Cells.Item(78, 12).Value = 20000
$ws.Cells.Item(78, 13).Value = 20000
$ws.Cells.Item(78, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(78, 15).Value = "Región del Maule"
$ws.Cells.Item(78, 16).Value = 800
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
